$wb = $excel.ActiveWorkbook

# --- Sheet "About": update the last-updated date in C1 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- Sheet "RAF-capacity": update hydrogen RAF values and view state ---
$wsCap = $wb.Worksheets.Item("RAF-capacity")
$wsCap.Range("B24").Value = 1
$wsCap.Range("B25").Value = 1

# Give column A a custom width on the RAF-capacity sheet
$wsCap.Columns.Item(1).ColumnWidth = 29.04296875

# Make RAF-capacity the active (selected) sheet/tab
$wsCap.Activate()
$wsCap.Select()

# Update the selection / scroll / zoom on RAF-capacity's view
$wsCap.Range("B25").Select()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.Zoom = 80
